$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Devices")

# Apply formatting first (copy format only) from representative cells,
# matching the style used by the real edit:
#   F1, G1  -> same style as header row cells (e.g. A7..I7, style index 8)
#   F2,F3,F4,G2 -> same style as B4 (style index 1)
$ws.Range("A7").Copy() | Out-Null
$ws.Range("F1:G1").PasteSpecial(-4122) | Out-Null

$ws.Range("B4").Copy() | Out-Null
$ws.Range("F2:G2").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("F3").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("F4").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Now set the cell values/content in the same order the strings were
# originally introduced so that the shared-strings table is rebuilt in the
# expected order: Loop, Built-in Loop-A, Built-in Loop-B, Column, Built-in Loop-C
$ws.Range("F1").Value = "Loop"
$ws.Range("F2").Value = "Built-in Loop-A"
$ws.Range("F3").Value = "Built-in Loop-B"
$ws.Range("G1").Value = "Column"
$ws.Range("F4").Value = "Built-in Loop-C"
$ws.Range("G2").Value = 2

# Update the active selection to F4, matching the edited workbook
$ws.Range("F4").Select() | Out-Null
